$wb = $excel.ActiveWorkbook

# --- Sheet "dbo#ExcelTest": remove the "Global Tolerance" row (row 3) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(3).Select()
$ws1.Rows.Item(3).Delete()

# --- Sheet "Assert": same fix (identical layout), then leave it as the active sheet/tab ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(3).Select()
$ws2.Rows.Item(3).Delete()

$ws2.Activate()
